$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the old last row (row 8), so it becomes row 11.
# This keeps the "last row" border style (s=9/10/11) on the true last row,
# matching how the sheet always formats the final data row distinctly.
$ws.Rows("8:10").Insert()

# The 3 freshly inserted rows come back blank/unstyled; copy the formatting
# (borders/alignment) from row 7 - the last "regular" row - onto them so they
# reuse the same existing style indices instead of Excel minting new ones.
$ws.Range("A7:N7").Copy()
$ws.Range("A8:N10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Helper: write literal text into a cell without Excels auto-conversion of
# date-shaped strings (e.g. "2020-03-11") into real dates.
function Set-TextCell($addr, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

# Row 2: DTHO-03-20-0032
$ws.Range("A2").Value = "DTHO"
$ws.Range("B2").Value = "Delito vs PExt"
Set-TextCell "C2" "2020-03-11"
$ws.Range("D2").Value = "Afectación a la planta exterior por corte de FON aérea. Localidad Arroyo seco"
$ws.Range("E2").Value = "Mayarí"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = "Fibra óptica-1"
$ws.Range("K2").Value = "No"
$ws.Range("L2").Value = "No"
$ws.Range("M2").Value = "No"
$ws.Range("N2").Value = "DTHO-03-20-0032"

# Row 3: DTHO-05-20-0046
$ws.Range("A3").Value = "DTHO"
$ws.Range("B3").Value = "Delito vs PExt"
Set-TextCell "C3" "2020-05-13"
$ws.Range("D3").Value = "Afectación a la planta exterior, por corte y sustracción de bajante telefónico. Calle Nicio García, No. 230,  entre 19 y 25"
$ws.Range("E3").Value = "Holguín"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 6.4
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = "Bajante telefónico-20"
$ws.Range("K3").Value = "No"
$ws.Range("L3").Value = "No"
$ws.Range("M3").Value = "No"
$ws.Range("N3").Value = "DTHO-05-20-0046"

# Row 4: DTHO-06-20-0059
$ws.Range("A4").Value = "DTHO"
$ws.Range("B4").Value = "Delito vs PExt"
Set-TextCell "C4" "2020-06-20"
$ws.Range("D4").Value = "Afectación a la planta exterior, por corte intensional de FO. Carretera vía San Andrés.Entrada a Purnio"
$ws.Range("E4").Value = "Holguín"
$ws.Range("F4").Value = 78.9
$ws.Range("G4").Value = 158.6
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = "Fibra óptica-1"
$ws.Range("K4").Value = "No"
$ws.Range("L4").Value = "No"
$ws.Range("M4").Value = "No"
$ws.Range("N4").Value = "DTHO-06-20-0059"

# Row 5: DTHO-07-20-0063
$ws.Range("A5").Value = "DTHO"
$ws.Range("B5").Value = "Delito vs PExt"
Set-TextCell "C5" "2020-07-10"
$ws.Range("D5").Value = "Afectación a la planta exterior, por corte y sustracción de bajante telefónico. Calle 26 de Julio,  No. 143"
$ws.Range("E5").Value = "Calixto García"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = "Bajante telefónico-48"
$ws.Range("K5").Value = "No"
$ws.Range("L5").Value = "No"
$ws.Range("M5").Value = "No"
$ws.Range("N5").Value = "DTHO-07-20-0063"

# Row 6: DTHO-07-20-0065
$ws.Range("A6").Value = "DTHO"
$ws.Range("B6").Value = "Delito vs PExt"
Set-TextCell "C6" "2020-07-14"
$ws.Range("D6").Value = "Afectación a la planta exterior, por corte y sustracción de bajante telefónico. Calle 13 de Marzo"
$ws.Range("E6").Value = "Calixto García"
$ws.Range("F6").Value = 50.3
$ws.Range("G6").Value = 257.6
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = "Bajante telefónico-200"
$ws.Range("K6").Value = "No"
$ws.Range("L6").Value = "No"
$ws.Range("M6").Value = "No"
$ws.Range("N6").Value = "DTHO-07-20-0065"

# Row 7: DTHO-07-20-0068
$ws.Range("A7").Value = "DTHO"
$ws.Range("B7").Value = "Delito vs PExt"
Set-TextCell "C7" "2020-07-26"
$ws.Range("D7").Value = "Afectación a la planta exterior, por Gabinete Distribución vandalizado. Carretera a San Germán. Rpto. 26 de Julio"
$ws.Range("E7").Value = "Holguín"
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 350
$ws.Range("I7").Value = "Gabinetes-1"
$ws.Range("K7").Value = "No"
$ws.Range("L7").Value = "No"
$ws.Range("M7").Value = "No"
$ws.Range("N7").Value = "DTHO-07-20-0068"

# Row 8: DTHO-08-20-0084
$ws.Range("A8").Value = "DTHO"
$ws.Range("B8").Value = "Delito vs PExt"
Set-TextCell "C8" "2020-08-28"
$ws.Range("D8").Value = "Afectación a la planta exterior por sustracción de cable de 400 pares. Carretera central Vía Las Tunas. Entre 36 y 38"
$ws.Range("E8").Value = "Holguín"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 394.23
$ws.Range("H8").Value = 146
$ws.Range("I8").Value = "Cable telefónico-60"
$ws.Range("K8").Value = "No"
$ws.Range("L8").Value = "No"
$ws.Range("M8").Value = "No"
$ws.Range("N8").Value = "DTHO-08-20-0084"

# Row 9: DTHO-09-20-0086
$ws.Range("A9").Value = "DTHO"
$ws.Range("B9").Value = "Delito vs PExt"
Set-TextCell "C9" "2020-09-03"
$ws.Range("D9").Value = "Afectación a la planta exterior por corte y sustracción de bajante telefónico. Localidad de San Germán. La Loma"
$ws.Range("E9").Value = "Urbano Noris"
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 6.52
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = "Bajante telefónico-410"
$ws.Range("K9").Value = "No"
$ws.Range("L9").Value = "No"
$ws.Range("M9").Value = "No"
$ws.Range("N9").Value = "DTHO-09-20-0086"

# Row 10: DTHO-11-20-0114
$ws.Range("A10").Value = "DTHO"
$ws.Range("B10").Value = "Delito vs PExt"
Set-TextCell "C10" "2020-11-12"
$ws.Range("D10").Value = "Afectación a la planta exterior por corte y sustracción de bajante telefónico. Vía las Tunas"
$ws.Range("E10").Value = "Holguín"
$ws.Range("F10").Value = 25
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "Bajante telefónico-50"
$ws.Range("K10").Value = "No"
$ws.Range("L10").Value = "No"
$ws.Range("M10").Value = "No"
$ws.Range("N10").Value = "DTHO-11-20-0114"

# Row 11: DTHO-11-20-0116
$ws.Range("A11").Value = "DTHO"
$ws.Range("B11").Value = "Delito vs PExt"
Set-TextCell "C11" "2020-11-24"
$ws.Range("D11").Value = "Afectación a la planta exterior por corte y sustracción de bajante telefónico. Carretera Central . Rpto. Oscar Lucero"
$ws.Range("E11").Value = "Holguín"
$ws.Range("F11").Value = 25
$ws.Range("G11").Value = 82.9
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = "Bajante telefónico-50"
$ws.Range("K11").Value = "No"
$ws.Range("L11").Value = "No"
$ws.Range("M11").Value = "No"
$ws.Range("N11").Value = "DTHO-11-20-0116"
